$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 49.666668
$ws.Range("I6").Value = 24.5
$ws.Range("K6").Value = 73.5
$ws.Range("M6").Value = 38.5

$ws.Range("H11").Value = 13.25
$ws.Range("I11").Value = 13.25
$ws.Range("K11").Value = 13.25
$ws.Range("M11").Value = 126.75

$ws.Range("H38").Value = 543
$ws.Range("I38").Value = 71.666664
$ws.Range("J38").Value = 1250
$ws.Range("K38").Value = 214.999992
$ws.Range("L38").Value = 3750
$ws.Range("M38").Value = 157.000008
$ws.Range("N38").Value = -4494

$ws.Range("H39").Value = 111.125
$ws.Range("I39").Value = 111.125
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 333.375
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -37.375
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 1000000
$ws.Range("J49").Value = 1000000
$ws.Range("L49").Value = 3000000
$ws.Range("N49").Value = -3000272

$ws.Range("H62").Value = 2652.6
$ws.Range("I62").Value = 2100
$ws.Range("J62").Value = 2790.75
$ws.Range("K62").Value = 2100
$ws.Range("L62").Value = 2790.75
$ws.Range("M62").Value = -1476
$ws.Range("N62").Value = -4038.75

$ws.Range("H65").Value = 2652.6
$ws.Range("I65").Value = 2100
$ws.Range("J65").Value = 2790.75
$ws.Range("K65").Value = 10500
$ws.Range("L65").Value = 13953.75
$ws.Range("M65").Value = -7380
$ws.Range("N65").Value = -20193.75

$ws.Range("H70").Value = 1752.7778
$ws.Range("I70").Value = 1591.6666
$ws.Range("J70").Value = 1833.3334
$ws.Range("K70").Value = 4774.9998
$ws.Range("L70").Value = 5500.0002
$ws.Range("M70").Value = -4504.9998
$ws.Range("N70").Value = -6040.0002

$ws.Range("H73").Value = 1752.7778
$ws.Range("I73").Value = 1591.6666
$ws.Range("J73").Value = 1833.3334
$ws.Range("K73").Value = 4774.9998
$ws.Range("L73").Value = 5500.0002
$ws.Range("M73").Value = -3838.9998
$ws.Range("N73").Value = -7372.0002

$ws.Range("H115").Value = 5666.3335
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

$ws.Range("H138").Value = 14029.566
$ws.Range("J138").Value = 14294.296
$ws.Range("L138").Value = 42882.888
$ws.Range("N138").Value = -53162.888

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1816.5
$ws.Range("J2").Value = 1547.25
$ws.Range("L2").Value = 1547.25
$ws.Range("N2").Value = -1773.25

$ws.Range("H45").Value = 1994.9286
$ws.Range("I45").Value = 1994.9286
$ws.Range("K45").Value = 1994.9286
$ws.Range("M45").Value = -1617.9286

$ws.Range("H61").Value = 2361.4
$ws.Range("I61").Value = 2361.4
$ws.Range("K61").Value = 2361.4
$ws.Range("M61").Value = -2149.4

$ws.Range("H74").Value = 2202
$ws.Range("I74").Value = 1442.4
$ws.Range("J74").Value = 6000
$ws.Range("K74").Value = 1442.4
$ws.Range("L74").Value = 6000
$ws.Range("M74").Value = -568.4000000000001
$ws.Range("N74").Value = -7748

$ws.Range("H77").Value = 2202
$ws.Range("I77").Value = 1442.4
$ws.Range("J77").Value = 6000
$ws.Range("K77").Value = 7212
$ws.Range("L77").Value = 30000
$ws.Range("M77").Value = -2844
$ws.Range("N77").Value = -38736

$ws.Range("H116").Value = 1816.5
$ws.Range("J116").Value = 1547.25
$ws.Range("L116").Value = 1547.25
$ws.Range("N116").Value = -6135.25

$ws.Range("H122").Value = 2516
$ws.Range("I122").Value = 2621.3333
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 7863.999899999999
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -5413.999899999999
$ws.Range("N122").Value = -11500

$ws.Range("H132").Value = 2709.6843
$ws.Range("I132").Value = 2250.5715
$ws.Range("K132").Value = 6751.7145
$ws.Range("M132").Value = -4221.7145

$ws.Range("H136").Value = 2361.4
$ws.Range("I136").Value = 2361.4
$ws.Range("K136").Value = 7084.200000000001
$ws.Range("M136").Value = -4534.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1816.5
$ws.Range("J3").Value = 1547.25
$ws.Range("L3").Value = 1547.25
$ws.Range("N3").Value = -1775.25

$ws.Range("H80").Value = 356
$ws.Range("J80").Value = 531.1429000000001
$ws.Range("L80").Value = 531.1429000000001
$ws.Range("N80").Value = -2527.1429

$ws.Range("H83").Value = 356
$ws.Range("J83").Value = 531.1429000000001
$ws.Range("L83").Value = 2655.7145
$ws.Range("N83").Value = -12639.7145

$ws.Range("H134").Value = 3087.625
$ws.Range("I134").Value = 3087.625
$ws.Range("K134").Value = 9262.875
$ws.Range("M134").Value = -6727.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 93325.37
$ws.Range("I16").Value = 113175.445
$ws.Range("K16").Value = 113175.445
$ws.Range("M16").Value = -112888.445

$ws.Range("H62").Value = 8962.25
$ws.Range("I62").Value = 1949.6666
$ws.Range("J62").Value = 30000
$ws.Range("K62").Value = 1949.6666
$ws.Range("L62").Value = 30000
$ws.Range("M62").Value = -1325.6666
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 8962.25
$ws.Range("I65").Value = 1949.6666
$ws.Range("J65").Value = 30000
$ws.Range("K65").Value = 9748.333000000001
$ws.Range("L65").Value = 150000
$ws.Range("M65").Value = -6628.333000000001
$ws.Range("N65").Value = -156240

$ws.Range("H105").Value = 2311.25
$ws.Range("J105").Value = 1299
$ws.Range("L105").Value = 1299
$ws.Range("N105").Value = -4793

$ws.Range("H107").Value = 92263.17999999999
$ws.Range("I107").Value = 126399.5
$ws.Range("K107").Value = 126399.5
$ws.Range("M107").Value = -124479.5

$ws.Range("H113").Value = 93325.37
$ws.Range("I113").Value = 113175.445
$ws.Range("K113").Value = 113175.445
$ws.Range("M113").Value = -111005.445

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 423660.75
$ws.Range("J4").Value = 626438.9
$ws.Range("L4").Value = 1879316.7
$ws.Range("N4").Value = -1879540.7

$ws.Range("H10").Value = 45
$ws.Range("I10").Value = 45
$ws.Range("K10").Value = 135
$ws.Range("M10").Value = 4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3013.3333
$ws.Range("I122").Value = 3013.3333
$ws.Range("K122").Value = 9039.999899999999
$ws.Range("M122").Value = -6589.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 8501560
$ws.Range("I61").Value = 6376665.5
$ws.Range("J61").Value = 12751350
$ws.Range("K61").Value = 6376665.5
$ws.Range("L61").Value = 12751350
$ws.Range("M61").Value = -6376463.5
$ws.Range("N61").Value = -12751754

$ws.Range("H113").Value = 8501560
$ws.Range("I113").Value = 6376665.5
$ws.Range("J113").Value = 12751350
$ws.Range("K113").Value = 6376665.5
$ws.Range("L113").Value = 12751350
$ws.Range("M113").Value = -6374495.5
$ws.Range("N113").Value = -12755690

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 3545.2727
$ws.Range("I132").Value = 3221.889
$ws.Range("K132").Value = 9665.667000000001
$ws.Range("M132").Value = -7135.667000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2775.55
$ws.Range("I132").Value = 2553.2942
$ws.Range("J132").Value = 4035
$ws.Range("K132").Value = 7659.882599999999
$ws.Range("L132").Value = 12105
$ws.Range("M132").Value = -5129.882599999999
$ws.Range("N132").Value = -17165
